$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The description for "cust_rating" (row 9, column B) had its number
# formatting re-worded from integers ("0 and 5") to decimals ("1.0 and 5.0").
$ws.Range("B9").Value = "Qualitative variable. Google rating of the restaurant if the number of customer reviews exceed 50. Possible values are between 1.0 and 5.0."
